$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H holds the tenant_id related header/cell (comment.tenant_id_lbl /
# model.tenant_id_lbl). Delete the entire column so remaining columns
# (update_usr_id_lbl, update_time_lbl) shift left from I/J to H/I.
$ws.Range("H1").EntireColumn.Delete()
